# DFII10.xlsx - append the two newest weekly observations to the
# "Quarterly" sheet (rows 96 and 97), matching the formatting already
# used by the rest of the A:B data columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quarterly")

# New observation_date / DFII10 rows
$ws.Range("A96").Value = 45966   # 2025-11-05
$ws.Range("B96").Value = 1.87
$ws.Range("A97").Value = 45973   # 2025-11-12
$ws.Range("B97").Value = 1.81

# Column A uses the custom yyyy-mm-dd date format (style seen on A2:A95);
# column B uses a plain 2-decimal number format (style seen on B2:B95).
$ws.Range("A96:A97").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B96:B97").NumberFormat = "0.00"

# Reflect the scrolled / selected state of the sheet after the append.
$ws.Range("B98").Select()
$excel.ActiveWindow.ScrollRow = 83
$excel.ActiveWindow.ScrollColumn = 1

# Best-effort: reflect the saved window geometry (may be a no-op in a
# headless/sandboxed host with no real screen window).
$excel.ActiveWindow.Left = 0
$excel.ActiveWindow.Top = 1470
$excel.ActiveWindow.Width = 23040
$excel.ActiveWindow.Height = 11670
